# Ally updates to triggers
# Adds cueMarker / targetMarker columns (E & F) with numeric trigger codes
# corresponding to each row's ConditionFile/TrialType grouping.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Headers
$ws.Range("E1").Value = "cueMarker"
$ws.Range("F1").Value = "targetMarker"

# Row -> (cueMarker, targetMarker) mapping based on trial groupings
$map = @{
    2  = @(1,6);  3  = @(1,6);  4  = @(1,6);  5  = @(1,6);
    6  = @(1,6);  7  = @(1,6);  8  = @(1,6);  9  = @(1,6);
    10 = @(2,7);  11 = @(2,7);  12 = @(2,7);  13 = @(2,7);
    14 = @(3,8);  15 = @(3,8);  16 = @(3,8);  17 = @(3,8);
    18 = @(3,8);  19 = @(3,8);  20 = @(3,8);  21 = @(3,8);
    22 = @(4,9);  23 = @(4,9);
    24 = @(5,10); 25 = @(5,10);
}

foreach ($r in $map.Keys) {
    $vals = $map[$r]
    $ws.Cells.Item($r, 5).Value = $vals[0]
    $ws.Cells.Item($r, 6).Value = $vals[1]
}

# Update the sheet view to match the new scroll/selection state
$ws.Range("F25").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
